$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=data(text) B=quantidade_atipica C=cliente D=id_venda(text) E=id_produto F=produto G=estoque_atualizado H=media_vendas I=desvio_padrao

# Row 2
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "2025-08-07"
$ws.Cells.Item(2,1).Style = "Normal"
$ws.Cells.Item(2,2).Value = 2
$ws.Cells.Item(2,3).Value = "BEMOL S/A"
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "393760"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = 14152
$ws.Cells.Item(2,6).Value = "HEADSET GAMER PLAYER PLUS LED 7 CORES 2M DRIVER 50MM PRETO LETRON"
$ws.Cells.Item(2,7).Value = -13
$ws.Cells.Item(2,8).Value = 1.08
$ws.Cells.Item(2,9).Value = 0.29

# Row 3
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2025-08-07"
$ws.Cells.Item(3,1).Style = "Normal"
$ws.Cells.Item(3,2).Value = 3
$ws.Cells.Item(3,3).Value = "BEMOL S/A"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "393791"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = 13965
$ws.Cells.Item(3,6).Value = "MINI VENTILADOR RECARREGAVEL E PORTATIL, COM PREGADOR MATERIAL ABS E COMPONENTES ELETRONICOS"
$ws.Cells.Item(3,7).Value = -26
$ws.Cells.Item(3,8).Value = 1.19
$ws.Cells.Item(3,9).Value = 0.49

# Row 4
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "2025-08-08"
$ws.Cells.Item(4,1).Style = "Normal"
$ws.Cells.Item(4,2).Value = 2
$ws.Cells.Item(4,3).Value = "BEMOL S/A"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "394429"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = 10130
$ws.Cells.Item(4,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(4,7).Value = -1370
$ws.Cells.Item(4,8).Value = 1.07
$ws.Cells.Item(4,9).Value = 0.3

# Row 5
$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "2025-08-11"
$ws.Cells.Item(5,1).Style = "Normal"
$ws.Cells.Item(5,2).Value = 2
$ws.Cells.Item(5,3).Value = "BEMOL S/A"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "396518"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = 10130
$ws.Cells.Item(5,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(5,7).Value = -1370
$ws.Cells.Item(5,8).Value = 1.07
$ws.Cells.Item(5,9).Value = 0.3

# Row 6
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "2025-08-11"
$ws.Cells.Item(6,1).Style = "Normal"
$ws.Cells.Item(6,2).Value = 2
$ws.Cells.Item(6,3).Value = "BEMOL S/A"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "396572"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = 10130
$ws.Cells.Item(6,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(6,7).Value = -1370
$ws.Cells.Item(6,8).Value = 1.07
$ws.Cells.Item(6,9).Value = 0.3

# Row 7
$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "2025-08-12"
$ws.Cells.Item(7,1).Style = "Normal"
$ws.Cells.Item(7,2).Value = 2
$ws.Cells.Item(7,3).Value = "BEMOL S/A"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "397270"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = 13544
$ws.Cells.Item(7,6).Value = "MOUSE SEM FIO 3 BOTOES 1000DPI COLOR FIT BRANCO 1709 R8"
$ws.Cells.Item(7,7).Value = 2
$ws.Cells.Item(7,8).Value = 1.07
$ws.Cells.Item(7,9).Value = 0.26

# Row 8
$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "2025-08-12"
$ws.Cells.Item(8,1).Style = "Normal"
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(8,3).Value = "BEMOL S/A"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "397270"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = 13244
$ws.Cells.Item(8,6).Value = "MOUSE SEM FIO 3 BOTOES 1000DPI COLOR FIT AZUL 1709 R8"
$ws.Cells.Item(8,7).Value = -15
$ws.Cells.Item(8,8).Value = 1.06
$ws.Cells.Item(8,9).Value = 0.25

# Row 9
$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = "2025-08-12"
$ws.Cells.Item(9,1).Style = "Normal"
$ws.Cells.Item(9,2).Value = 2
$ws.Cells.Item(9,3).Value = "BEMOL S/A"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "397297"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = 12945
$ws.Cells.Item(9,6).Value = "FONE DE OUVIDO SEM FIO BT BASIKE FON-9856"
$ws.Cells.Item(9,7).Value = -107
$ws.Cells.Item(9,8).Value = 1.03
$ws.Cells.Item(9,9).Value = 0.18

# Row 10
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "2025-08-12"
$ws.Cells.Item(10,1).Style = "Normal"
$ws.Cells.Item(10,2).Value = 2
$ws.Cells.Item(10,3).Value = "BEMOL S/A"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "397345"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = 396985
$ws.Cells.Item(10,6).Value = "Kit Smartwatch Inova Com Pulseira Respiravel Preto Fone E Carregador"
$ws.Cells.Item(10,7).Value = -59
$ws.Cells.Item(10,8).Value = 1.02
$ws.Cells.Item(10,9).Value = 0.13

# Row 11
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "2025-08-13"
$ws.Cells.Item(11,1).Style = "Normal"
$ws.Cells.Item(11,2).Value = 2
$ws.Cells.Item(11,3).Value = "BEMOL S/A"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "398131"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = 10130
$ws.Cells.Item(11,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(11,7).Value = -1370
$ws.Cells.Item(11,8).Value = 1.07
$ws.Cells.Item(11,9).Value = 0.3

# Row 12
$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = "2025-08-13"
$ws.Cells.Item(12,1).Style = "Normal"
$ws.Cells.Item(12,2).Value = 2
$ws.Cells.Item(12,3).Value = "BEMOL S/A"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "398157"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = 11436
$ws.Cells.Item(12,6).Value = "FONE HEADSET FONE DE OUVIDO PEI-P9 MUSIC POWER"
$ws.Cells.Item(12,7).Value = -2
$ws.Cells.Item(12,8).Value = 1.02
$ws.Cells.Item(12,9).Value = 0.14

# Row 13
$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = "2025-08-18"
$ws.Cells.Item(13,1).Style = "Normal"
$ws.Cells.Item(13,2).Value = 2
$ws.Cells.Item(13,3).Value = "BEMOL S/A"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "400531"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = 418713
$ws.Cells.Item(13,6).Value = "SAPATEIRA 12 PARES FERRO + PPATE 2,5KG, TAMANHO 50X18,5X55CM"
$ws.Cells.Item(13,7).Value = -24
$ws.Cells.Item(13,8).Value = 1.09
$ws.Cells.Item(13,9).Value = 0.29

# Row 14
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2025-08-18"
$ws.Cells.Item(14,1).Style = "Normal"
$ws.Cells.Item(14,2).Value = 2
$ws.Cells.Item(14,3).Value = "BEMOL S/A"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "400543"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = 418713
$ws.Cells.Item(14,6).Value = "SAPATEIRA 12 PARES FERRO + PPATE 2,5KG, TAMANHO 50X18,5X55CM"
$ws.Cells.Item(14,7).Value = -24
$ws.Cells.Item(14,8).Value = 1.09
$ws.Cells.Item(14,9).Value = 0.29

# Row 15
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2025-08-18"
$ws.Cells.Item(15,1).Style = "Normal"
$ws.Cells.Item(15,2).Value = 2
$ws.Cells.Item(15,3).Value = "BEMOL S/A"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "400578"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = 10525
$ws.Cells.Item(15,6).Value = "BOMBA AUTOMATICA PARA GALAO DE AGUA RECARREGAVEL USB"
$ws.Cells.Item(15,7).Value = -133
$ws.Cells.Item(15,8).Value = 1.04
$ws.Cells.Item(15,9).Value = 0.2

# Row 16
$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2025-08-18"
$ws.Cells.Item(16,1).Style = "Normal"
$ws.Cells.Item(16,2).Value = 2
$ws.Cells.Item(16,3).Value = "BEMOL S/A"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "400733"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = 14280
$ws.Cells.Item(16,6).Value = "SUPORTE DE MESA DOBRÁVEL ROTAÇÃO 360° HMASTON"
$ws.Cells.Item(16,7).Value = -5
$ws.Cells.Item(16,8).Value = 1.08
$ws.Cells.Item(16,9).Value = 0.28

# Row 17
$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2025-08-18"
$ws.Cells.Item(17,1).Style = "Normal"
$ws.Cells.Item(17,2).Value = 2
$ws.Cells.Item(17,3).Value = "BEMOL S/A"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "400752"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = 13198
$ws.Cells.Item(17,6).Value = "ESCOVA DE LIMPEZA ELETRICA RECARREGAVEL ALIMENTACAO BATERIA 18650 LI-ION TAMANHO 20X7,5X6,2CM"
$ws.Cells.Item(17,7).Value = -85
$ws.Cells.Item(17,8).Value = 1.05
$ws.Cells.Item(17,9).Value = 0.21

# Row 18
$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2025-08-19"
$ws.Cells.Item(18,1).Style = "Normal"
$ws.Cells.Item(18,2).Value = 3
$ws.Cells.Item(18,3).Value = "BEMOL S/A"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "401315"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = 13854
$ws.Cells.Item(18,6).Value = "CANETA APAGAVEL GEL 0.7 AZUL/PRETA - CORES SORTIDAS JOCAR OFFICE"
$ws.Cells.Item(18,7).Value = 0
$ws.Cells.Item(18,8).Value = 1.24
$ws.Cells.Item(18,9).Value = 0.56

# Row 19
$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = "2025-08-19"
$ws.Cells.Item(19,1).Style = "Normal"
$ws.Cells.Item(19,2).Value = 2
$ws.Cells.Item(19,3).Value = "BEMOL S/A"
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "401319"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = 12945
$ws.Cells.Item(19,6).Value = "FONE DE OUVIDO SEM FIO BT BASIKE FON-9856"
$ws.Cells.Item(19,7).Value = -107
$ws.Cells.Item(19,8).Value = 1.03
$ws.Cells.Item(19,9).Value = 0.18

# Row 20
$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = "2025-08-19"
$ws.Cells.Item(20,1).Style = "Normal"
$ws.Cells.Item(20,2).Value = 2
$ws.Cells.Item(20,3).Value = "BEMOL S/A"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "401332"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = 10525
$ws.Cells.Item(20,6).Value = "BOMBA AUTOMATICA PARA GALAO DE AGUA RECARREGAVEL USB"
$ws.Cells.Item(20,7).Value = -133
$ws.Cells.Item(20,8).Value = 1.04
$ws.Cells.Item(20,9).Value = 0.2

# Row 21
$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = "2025-08-19"
$ws.Cells.Item(21,1).Style = "Normal"
$ws.Cells.Item(21,2).Value = 2
$ws.Cells.Item(21,3).Value = "BEMOL S/A"
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "401336"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = 13977
$ws.Cells.Item(21,6).Value = "SAPATEIRA MATERIAL PP, FERRO E TNT CAPACIDADE18 PARES, SUPORTA ATE 15KG"
$ws.Cells.Item(21,7).Value = -38
$ws.Cells.Item(21,8).Value = 1.03
$ws.Cells.Item(21,9).Value = 0.16

# Row 22
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = "2025-08-19"
$ws.Cells.Item(22,1).Style = "Normal"
$ws.Cells.Item(22,2).Value = 2
$ws.Cells.Item(22,3).Value = "BEMOL S/A"
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "401337"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = 10130
$ws.Cells.Item(22,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(22,7).Value = -1370
$ws.Cells.Item(22,8).Value = 1.07
$ws.Cells.Item(22,9).Value = 0.3

# Row 23
$ws.Cells.Item(23,1).NumberFormat = "@"
$ws.Cells.Item(23,1).Value = "2025-08-19"
$ws.Cells.Item(23,1).Style = "Normal"
$ws.Cells.Item(23,2).Value = 2
$ws.Cells.Item(23,3).Value = "BEMOL S/A"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "401348"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = 10130
$ws.Cells.Item(23,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(23,7).Value = -1370
$ws.Cells.Item(23,8).Value = 1.07
$ws.Cells.Item(23,9).Value = 0.3

# Row 24
$ws.Cells.Item(24,1).NumberFormat = "@"
$ws.Cells.Item(24,1).Value = "2025-08-19"
$ws.Cells.Item(24,1).Style = "Normal"
$ws.Cells.Item(24,2).Value = 3
$ws.Cells.Item(24,3).Value = "BEMOL S/A"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "401363"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = 10130
$ws.Cells.Item(24,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(24,7).Value = -1370
$ws.Cells.Item(24,8).Value = 1.07
$ws.Cells.Item(24,9).Value = 0.3

# Row 25
$ws.Cells.Item(25,1).NumberFormat = "@"
$ws.Cells.Item(25,1).Value = "2025-08-19"
$ws.Cells.Item(25,1).Style = "Normal"
$ws.Cells.Item(25,2).Value = 2
$ws.Cells.Item(25,3).Value = "BEMOL S/A"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "401365"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = 10130
$ws.Cells.Item(25,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(25,7).Value = -1370
$ws.Cells.Item(25,8).Value = 1.07
$ws.Cells.Item(25,9).Value = 0.3

# Row 26
$ws.Cells.Item(26,1).NumberFormat = "@"
$ws.Cells.Item(26,1).Value = "2025-08-19"
$ws.Cells.Item(26,1).Style = "Normal"
$ws.Cells.Item(26,2).Value = 2
$ws.Cells.Item(26,3).Value = "BEMOL S/A"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "401953"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = 48646
$ws.Cells.Item(26,6).Value = "COLHER MEDIDORA COM BALANÇA DIGITAL PRETO"
$ws.Cells.Item(26,7).Value = -13
$ws.Cells.Item(26,8).Value = 1.08
$ws.Cells.Item(26,9).Value = 0.29

# Row 27
$ws.Cells.Item(27,1).NumberFormat = "@"
$ws.Cells.Item(27,1).Value = "2025-08-19"
$ws.Cells.Item(27,1).Style = "Normal"
$ws.Cells.Item(27,2).Value = 2
$ws.Cells.Item(27,3).Value = "BEMOL S/A"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "401965"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = 10130
$ws.Cells.Item(27,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(27,7).Value = -1370
$ws.Cells.Item(27,8).Value = 1.07
$ws.Cells.Item(27,9).Value = 0.3

# Row 28
$ws.Cells.Item(28,1).NumberFormat = "@"
$ws.Cells.Item(28,1).Value = "2025-08-19"
$ws.Cells.Item(28,1).Style = "Normal"
$ws.Cells.Item(28,2).Value = 3
$ws.Cells.Item(28,3).Value = "BEMOL S/A"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "401969"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = 12651
$ws.Cells.Item(28,6).Value = "CARREGADOR PARA NOTEBBOK AGOLD 120W"
$ws.Cells.Item(28,7).Value = -118
$ws.Cells.Item(28,8).Value = 1.02
$ws.Cells.Item(28,9).Value = 0.22

# Row 29
$ws.Cells.Item(29,1).NumberFormat = "@"
$ws.Cells.Item(29,1).Value = "2025-08-21"
$ws.Cells.Item(29,1).Style = "Normal"
$ws.Cells.Item(29,2).Value = 3
$ws.Cells.Item(29,3).Value = "BEMOL S/A"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "402635"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = 10130
$ws.Cells.Item(29,6).Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Cells.Item(29,7).Value = -1370
$ws.Cells.Item(29,8).Value = 1.07
$ws.Cells.Item(29,9).Value = 0.3

# Row 30
$ws.Cells.Item(30,1).NumberFormat = "@"
$ws.Cells.Item(30,1).Value = "2025-08-21"
$ws.Cells.Item(30,1).Style = "Normal"
$ws.Cells.Item(30,2).Value = 2
$ws.Cells.Item(30,3).Value = "BEMOL S/A"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "402640"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = 12886
$ws.Cells.Item(30,6).Value = "MASSAGEADOR MUSCULAR INOVA JMQ-12216"
$ws.Cells.Item(30,7).Value = -154
$ws.Cells.Item(30,8).Value = 1.01
$ws.Cells.Item(30,9).Value = 0.12
